# Insert a new weekly price observation as row 444 in the data table,
# shifting the existing rows 444-547 down to 445-548.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 444 (this shifts rows 444..547 -> 445..548
# and copies formatting, e.g. the date number format, from the row above).
$ws.Rows.Item(444).Insert()

# Populate the newly inserted row 444 with the new observation.
$ws.Cells.Item(444, 1).Value  = 4
$ws.Cells.Item(444, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(444, 3).Value  = "Los Lagos"
$ws.Cells.Item(444, 4).Value  = 45244
$ws.Cells.Item(444, 5).Value  = 10
$ws.Cells.Item(444, 6).Value  = 100112037
$ws.Cells.Item(444, 7).Value  = "Cebollín"
$ws.Cells.Item(444, 8).Value  = "Sin especificar"
$ws.Cells.Item(444, 9).Value  = "Primera"
$ws.Cells.Item(444, 10).Value = 180
$ws.Cells.Item(444, 11).Value = 6500
$ws.Cells.Item(444, 12).Value = 6500
$ws.Cells.Item(444, 13).Value = 6500
$ws.Cells.Item(444, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(444, 15).Value = "Región Metropolitana"
$ws.Cells.Item(444, 16).Value = 181
$ws.Cells.Item(444, 17).Value = 36
$ws.Cells.Item(444, 18).Value = "Hortaliza"
